$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 2.44
$ws.Range("I2").Value = 2.46
$ws.Range("L2").Value = 1.34
$ws.Range("P2").Value = 1.91
$ws.Range("R2").Value = 1.36
$ws.Range("V2").Value = 1.68
$ws.Range("X2").Value = 17
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 16
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 12.5
$ws.Range("AC2").Value = 8.199999999999999
$ws.Range("AE2").Value = 27
$ws.Range("AI2").Value = 42
$ws.Range("AJ2").Value = 1000
$ws.Range("AM2").Value = 120
$ws.Range("AN2").Value = 36
$ws.Range("AO2").Value = 26
# Row 3
$ws.Range("F3").Value = 1.04
$ws.Range("G3").Value = 1000
$ws.Range("H3").Value = 1.04
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 1.01
$ws.Range("K3").Value = 980
$ws.Range("Q3").Value = 1.76
# Row 4
$ws.Range("F4").Value = 1.79
$ws.Range("G4").Value = 2.14
$ws.Range("H4").Value = 3.8
$ws.Range("I4").Value = 6.6
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 3.65
$ws.Range("P4").Value = 1.47
$ws.Range("Q4").Value = 2.4
# Row 5
$ws.Range("F5").Value = 3.15
$ws.Range("G5").Value = 5.7
$ws.Range("H5").Value = 2.12
$ws.Range("I5").Value = 2.68
$ws.Range("J5").Value = 2.46
$ws.Range("K5").Value = 3.55
$ws.Range("P5").Value = 1.36
$ws.Range("Q5").Value = 2.86
# Row 6
$ws.Range("F6").Value = 1.92
$ws.Range("H6").Value = 3.55
$ws.Range("J6").Value = 3.1
$ws.Range("K6").Value = 5.4
$ws.Range("P6").Value = 1.77
# Row 7
$ws.Range("F7").Value = 1.95
$ws.Range("G7").Value = 2.42
$ws.Range("H7").Value = 3.35
$ws.Range("I7").Value = 6.2
$ws.Range("J7").Value = 2.64
$ws.Range("K7").Value = 3.85
# Row 8
$ws.Range("F8").Value = 1.25
$ws.Range("G8").Value = 1000
$ws.Range("H8").Value = 1.25
$ws.Range("I8").Value = 1000
$ws.Range("J8").Value = 1.01
$ws.Range("K8").Value = 1000
$ws.Range("P8").Value = 2.4
# Row 9
$ws.Range("N9").Value = 5
$ws.Range("T9").Value = 1.62
$ws.Range("AJ9").Value = 65
# Row 10
$ws.Range("G10").Value = 3.65
$ws.Range("J10").Value = 2.94
$ws.Range("K10").Value = 2.98
$ws.Range("AG10").Value = 17.5
